$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (LinearRegression) - C2, D2 change
$ws.Range("C2").Value = -7.182238711347325
$ws.Range("D2").Value = -7.182238711347325

# Row 3 (RandomForestRegressor) - B3, C3, D3 change
$ws.Range("B3").Value = 0.7565973768040211
$ws.Range("C3").Value = 0.7410516786700642
$ws.Range("D3").Value = 0.04270157478783974

# Row 4 - label change GradientBoostingRegressor -> DecisionTreeRegressor, values change
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.862766372495713
$ws.Range("C4").Value = 0.8647582932951939
$ws.Range("D4").Value = -1.445468350413504

# Row 5 - label change AdaBoostRegressor -> MLPRegressor, values change
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.6797398819107157
$ws.Range("C5").Value = 0.5448335826361347
$ws.Range("D5").Value = -3.316651931249464
